$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 7 - Kommunikation med Anders
$ws.Range("A7").Value = "Kommunikation med Anders"
$ws.Range("C7").Value = 43963
$ws.Range("D7").Value = 0.354166666666667
$ws.Range("E7").Value = 0.395833333333333

# Row 8 - Grundopsætning af rapport
$ws.Range("A8").Value = "Grundopsætning af rapport"
$ws.Range("C8").Value = 43963
$ws.Range("D8").Value = 0.395833333333333
$ws.Range("E8").Value = 0.458333333333333

# Row 9 - Metrikker – bilag
$ws.Range("A9").Value = "Metrikker – bilag"
$ws.Range("C9").Value = 43963
$ws.Range("D9").Value = 0.458333333333333
$ws.Range("E9").Value = 0.541666666666667

# Row 10 - Review risikoanalyse
$ws.Range("A10").Value = "Review risikoanalyse"
$ws.Range("C10").Value = 43963
$ws.Range("D10").Value = 0.541666666666667
$ws.Range("E10").Value = 0.583333333333333

# Row 11 - Projektplan
$ws.Range("A11").Value = "Projektplan"
$ws.Range("C11").Value = 43963
$ws.Range("D11").Value = 0.583333333333333
$ws.Range("E11").Value = 0.677083333333333

# Row 12 - Metrikker – rapport
$ws.Range("A12").Value = "Metrikker – rapport"
$ws.Range("C12").Value = 43963
$ws.Range("D12").Value = 0.677083333333333
$ws.Range("E12").Value = 0.71875

# Move the active selection to C13, matching the author's last edit position
$ws.Range("C13").Select() | Out-Null
